# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Logical changes being applied to the "Estado de Cuenta" sheet:
#   - A new "Periodo Mora" row is added for period 2508 (a new row 19).
#   - The existing three period rows are re-ordered from descending
#     (2507, 2506, 2505) to ascending (2505, 2506, 2507) so the table now
#     reads 2505, 2506, 2507, 2508 top to bottom.
#   - The "VALOR MORA" total (E11) and "Cant. Periodos" count (F13) are
#     updated to reflect the new period.
#   - The signature block (line + labels) that used to sit on rows 23-24
#     is pushed down to rows 24-25 because of the inserted row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room for the new period row -------------------------------------
# Insert a blank row at 19; everything from row 19 down (the signature
# block on rows 23/24) shifts down by one, landing on rows 24/25.
$ws.Rows("19").Insert()

# --- Formatting ---------------------------------------------------------
# Row 18 currently still carries the special "bottom of table" border
# (it used to be the last data row). That look now belongs to the new
# last row (19), so copy it down first...
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)

# ...then restyle row 18 as a normal/interior row, matching rows 16-17.
$ws.Range("B16:J16").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)

# --- Content: add the new 2508 period row (row 19) -----------------------
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1044907855"
$ws.Range("D19").Value = "KENEL JOSE AURELA MARTINEZ"
$ws.Range("E19").Value = "2508"
$ws.Range("F19").Value = 62663
$ws.Range("G19").Value = 1566560

# --- Content: re-order the existing periods to ascending ------------------
$ws.Range("E16").Value = "2505"
$ws.Range("E17").Value = "2506"
$ws.Range("E18").Value = "2507"

# --- Update summary figures ------------------------------------------------
$ws.Range("E11").Value = 250652
$ws.Range("F13").Value = 4
